# Updates the "cryptos" price table (columns B-E) to the latest scrape.
# Numeric-looking text values (e.g. "582.75") must be written as literal
# text -- the source sheet stores Coin/Link/Price/Volume as strings, not
# numbers -- so Set-TextValue forces text via NumberFormat "@" and then
# clears that temporary format again so no stray cell style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -match "^[+-]?\d+(\.\d+)?$") {
        # Looks like a plain number -- force text so Excel does not
        # reinterpret it (and restore the default style afterwards).
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

# Row 2
Set-TextValue 2 4 "63.168.94"
Set-TextValue 2 5 "  +0.43%  "
# Row 3
Set-TextValue 3 4 "2.572.16"
Set-TextValue 3 5 "  +1.66%  "
# Row 4
Set-TextValue 4 5 "  +0.01%  "
# Row 5
Set-TextValue 5 4 "582.75"
Set-TextValue 5 5 "  +2.60%  "
# Row 6
Set-TextValue 6 4 "148.43"
Set-TextValue 6 5 "  +0.51%  "
# Row 7
Set-TextValue 7 4 "0.999"
Set-TextValue 7 5 "  +0.01%  "
# Row 8
Set-TextValue 8 5 "  +0.67%  "
# Row 9
Set-TextValue 9 4 "0.108"
Set-TextValue 9 5 "  +2.40%  "
# Row 10
Set-TextValue 10 4 "5.61"
Set-TextValue 10 5 "  -0.64%  "
# Row 11
Set-TextValue 11 5 "  +0.20%  "
# Row 12
Set-TextValue 12 5 "  +0.60%  "
# Row 13
Set-TextValue 13 4 "27.70"
Set-TextValue 13 5 "  +0.01%  "
# Row 14
Set-TextValue 14 4 "3.031.40"
Set-TextValue 14 5 "  +1.68%  "
# Row 15
Set-TextValue 15 4 "63.101.88"
Set-TextValue 15 5 "  +0.44%  "
# Row 16
Set-TextValue 16 5 "  +2.84%  "
# Row 17
Set-TextValue 17 4 "2.568.21"
Set-TextValue 17 5 "  +2.00%  "
# Row 18
Set-TextValue 18 4 "11.43"
Set-TextValue 18 5 "  -1.18%  "
# Row 19
Set-TextValue 19 4 "341.04"
Set-TextValue 19 5 "  +1.66%  "
# Row 20
Set-TextValue 20 5 "  +1.93%  "
# Row 21
Set-TextValue 21 4 "6.86"
Set-TextValue 21 5 "  +1.54%  "
# Row 22
Set-TextValue 22 5 "  +0.05%  "
# Row 23
Set-TextValue 23 4 "65.95"
Set-TextValue 23 5 "  +0.58%  "
# Row 24
Set-TextValue 24 4 "2.682.57"
Set-TextValue 24 5 "  +0.83%  "
# Row 25
Set-TextValue 25 5 "  +4.06%  "
# Row 26
Set-TextValue 26 5 "  +1.09%  "
# Row 27
Set-TextValue 27 2 "Aptos"
Set-TextValue 27 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 27 4 "8.03"
Set-TextValue 27 5 "  +10.90%  "
# Row 28
Set-TextValue 28 2 "InternetComputer(DFINITY)"
Set-TextValue 28 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 28 4 "8.51"
Set-TextValue 28 5 "  +1.96%  "
# Row 29
Set-TextValue 29 2 "Binance-PegBSC-USD"
Set-TextValue 29 3 "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue 29 4 "1.00"
Set-TextValue 29 5 "  +0.04%  "
# Row 30
Set-TextValue 30 2 "SuiNetwork"
Set-TextValue 30 3 "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue 30 4 "1.48"
Set-TextValue 30 5 "  -2.55%  "
# Row 31
Set-TextValue 31 5 "  +6.30%  "
# Row 32
Set-TextValue 32 4 "0.0₃0827"
Set-TextValue 32 5 "  +1.87%  "
# Row 33
Set-TextValue 33 4 "177.50"
Set-TextValue 33 5 "  -0.02%  "
# Row 34
Set-TextValue 34 4 "438.54"
Set-TextValue 34 5 "  +6.50%  "
# Row 35
Set-TextValue 35 4 "1.61"
Set-TextValue 35 5 "  +1.48%  "
# Row 36
Set-TextValue 36 5 "  +1.80%  "
# Row 37
Set-TextValue 37 4 "19.32"
Set-TextValue 37 5 "  +2.69%  "
# Row 38
Set-TextValue 38 4 "4.53"
Set-TextValue 38 5 "  +3.27%  "
# Row 39
Set-TextValue 39 5 "  +0.04%  "
# Row 40
Set-TextValue 40 4 "1.75"
Set-TextValue 40 5 "  +1.08%  "
# Row 41
Set-TextValue 41 5 "  +0.12%  "
# Row 42
Set-TextValue 42 4 "152.34"
Set-TextValue 42 5 "  +0.33%  "
# Row 43
Set-TextValue 43 4 "3.82"
Set-TextValue 43 5 "  +1.89%  "
# Row 44
Set-TextValue 44 4 "21.17"
Set-TextValue 44 5 "  +2.41%  "
# Row 45
Set-TextValue 45 5 "  +6.15%  "
# Row 46
Set-TextValue 46 4 "0.608"
Set-TextValue 46 5 "  +0.71%  "
# Row 47
Set-TextValue 47 4 "0.0975"
Set-TextValue 47 5 "  +0.96%  "
# Row 48
Set-TextValue 48 4 "0.0243"
Set-TextValue 48 5 "  +2.19%  "
# Row 49
Set-TextValue 49 4 "18.47"
Set-TextValue 49 5 "  +0.94%  "
# Row 50
Set-TextValue 50 4 "1.74"
Set-TextValue 50 5 "  -2.48%  "
# Row 51
Set-TextValue 51 4 "11.38"
Set-TextValue 51 5 "  -0.14%  "
